$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 8.69
$ws.Range("E2").Value = 59
$ws.Range("F2").Value = 11.59
$ws.Range("H2").Value = 50
$ws.Range("K2").Value = 57.2
$ws.Range("N2").Value = 50.60178744571824

# Row 3 updates
$ws.Range("D3").Value = 15.35
$ws.Range("E3").Value = 57.2
$ws.Range("F3").Value = 6.34
$ws.Range("K3").Value = 57.2
$ws.Range("N3").Value = 50.60178744571824
